$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = "repo-migration"
$ws.Range("C2").Value = "almatasks"
$ws.Range("D2").Value = "anilgoudasb06"
$ws.Range("E2").Value = "REP"

# Update row 3
$ws.Range("B3").Value = "repo-migration"
$ws.Range("C3").Value = "app-n-pak"
$ws.Range("D3").Value = "anilgoudasb06"
$ws.Range("E3").Value = "REP"

# Add new row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "repo-migration"
$ws.Range("C4").Value = "casa-build-utils"
$ws.Range("D4").Value = "anilgoudasb06"
$ws.Range("E4").Value = "REP"

# Add new row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "repo-migration"
$ws.Range("C5").Value = "casa6"
$ws.Range("D5").Value = "anilgoudasb06"
$ws.Range("E5").Value = "REP"

# Add new row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "repo-migration"
$ws.Range("C6").Value = "casashell"
$ws.Range("D6").Value = "anilgoudasb06"
$ws.Range("E6").Value = "REP"

# Move selection to B5 as in the final document
$ws.Range("B5").Select()
